$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in B2:E13 to the nearest integer, writing
# them back as plain integer data (matching how Ontpl_/Pot_ files are
# now persisted to disk).
for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val) {
            $cell.Value2 = $excel.WorksheetFunction.Round([double]$val, 0)
        }
    }
}
